# Auto-generated Excel COM-interop script
# Applies numeric cell updates (market price refresh) across all 8 sheets
# as described by the Omega_Profits.xlsx diff (scheduled runner data refresh).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 593.4167
$ws.Range("I28").Value = 593.4167
$ws.Range("K28").Value = 593.4167
$ws.Range("M28").Value = -108.4167
$ws.Range("H80").Value = 822.1429000000001
$ws.Range("I80").Value = 405.8
$ws.Range("J80").Value = 1053.4445
$ws.Range("K80").Value = 1217.4
$ws.Range("L80").Value = 3160.3335
$ws.Range("M80").Value = -219.4000000000001
$ws.Range("N80").Value = -5156.333500000001
$ws.Range("H83").Value = 822.1429000000001
$ws.Range("I83").Value = 405.8
$ws.Range("J83").Value = 1053.4445
$ws.Range("K83").Value = 3652.2
$ws.Range("L83").Value = 9481.0005
$ws.Range("M83").Value = 1339.8
$ws.Range("N83").Value = -19465.0005
$ws.Range("H96").Value = 597.7143
$ws.Range("I96").Value = 307.375
$ws.Range("J96").Value = 984.8333
$ws.Range("K96").Value = 922.125
$ws.Range("L96").Value = 2954.4999
$ws.Range("M96").Value = 450.875
$ws.Range("N96").Value = -5700.4999
$ws.Range("H100").Value = 4371.6665
$ws.Range("J100").Value = 5947.6665
$ws.Range("L100").Value = 5947.6665
$ws.Range("N100").Value = -7029.6665
$ws.Range("H107").Value = 1493.4286
$ws.Range("I107").Value = 1543.5385
$ws.Range("J107").Value = 842
$ws.Range("K107").Value = 1543.5385
$ws.Range("L107").Value = 842
$ws.Range("M107").Value = 376.4614999999999
$ws.Range("N107").Value = -4682
$ws.Range("H113").Value = 2462.5
$ws.Range("I113").Value = 2666.6667
$ws.Range("K113").Value = 2666.6667
$ws.Range("M113").Value = 587.3332999999998
$ws.Range("H115").Value = 779.75
$ws.Range("I115").Value = 779.75
$ws.Range("K115").Value = 2339.25
$ws.Range("M115").Value = -772.25
$ws.Range("H137").Value = 4281.7827
$ws.Range("J137").Value = 4530.5
$ws.Range("L137").Value = 13591.5
$ws.Range("N137").Value = -18691.5
$ws.Range("H138").Value = 9381.621999999999
$ws.Range("J138").Value = 9890.3125
$ws.Range("L138").Value = 29670.9375
$ws.Range("N138").Value = -39950.9375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 54488.11
$ws.Range("I32").Value = 90000
$ws.Range("J32").Value = 50049.125
$ws.Range("K32").Value = 90000
$ws.Range("L32").Value = 50049.125
$ws.Range("M32").Value = -89713
$ws.Range("N32").Value = -50623.125
$ws.Range("H37").Value = 18222.6
$ws.Range("J37").Value = 17778.25
$ws.Range("L37").Value = 17778.25
$ws.Range("N37").Value = -18324.25
$ws.Range("H44").Value = 49999
$ws.Range("J44").Value = 49999
$ws.Range("L44").Value = 49999
$ws.Range("N44").Value = -50975
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H102").Value = 2278.6667
$ws.Range("I102").Value = 2122.182
$ws.Range("K102").Value = 2122.182
$ws.Range("M102").Value = -500.1819999999998
$ws.Range("H110").Value = 1780.625
$ws.Range("I110").Value = 1780.625
$ws.Range("K110").Value = 1780.625
$ws.Range("M110").Value = 264.375
$ws.Range("H122").Value = 12584
$ws.Range("I122").Value = 12867.625
$ws.Range("K122").Value = 38602.875
$ws.Range("M122").Value = -36152.875
$ws.Range("H132").Value = 6132.8
$ws.Range("I132").Value = 6274.4707
$ws.Range("J132").Value = 5330
$ws.Range("K132").Value = 18823.4121
$ws.Range("L132").Value = 15990
$ws.Range("M132").Value = -16293.4121
$ws.Range("N132").Value = -21050

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1136.3
$ws.Range("I107").Value = 1107.2632
$ws.Range("K107").Value = 1107.2632
$ws.Range("M107").Value = 812.7367999999999
$ws.Range("H130").Value = 88333.336
$ws.Range("J130").Value = 88333.336
$ws.Range("L130").Value = 88333.336
$ws.Range("N130").Value = -98373.336

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2007.1428
$ws.Range("I134").Value = 2061.8235
$ws.Range("J134").Value = 1774.75
$ws.Range("K134").Value = 6185.470499999999
$ws.Range("L134").Value = 5324.25
$ws.Range("M134").Value = -3650.470499999999
$ws.Range("N134").Value = -10394.25
$ws.Range("H141").Value = 388618
$ws.Range("J141").Value = 388618
$ws.Range("L141").Value = 388618
$ws.Range("N141").Value = -398978

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1747.5238
$ws.Range("I5").Value = 875
$ws.Range("J5").Value = 2096.5334
$ws.Range("K5").Value = 2625
$ws.Range("L5").Value = 6289.600199999999
$ws.Range("M5").Value = -2513
$ws.Range("N5").Value = -6513.600199999999
$ws.Range("H98").Value = 857.2222
$ws.Range("I98").Value = 653.25
$ws.Range("J98").Value = 1020.4
$ws.Range("K98").Value = 1959.75
$ws.Range("L98").Value = 3061.2
$ws.Range("M98").Value = -461.75
$ws.Range("N98").Value = -6057.2
$ws.Range("H132").Value = 2660.375
$ws.Range("I132").Value = 1321
$ws.Range("K132").Value = 11889
$ws.Range("M132").Value = -9359
$ws.Range("H135").Value = 1747.5238
$ws.Range("I135").Value = 875
$ws.Range("J135").Value = 2096.5334
$ws.Range("K135").Value = 7875
$ws.Range("L135").Value = 18868.8006
$ws.Range("M135").Value = -5340
$ws.Range("N135").Value = -23938.8006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 641
$ws.Range("I97").Value = 635.63635
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 635.63635
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -139.63635
$ws.Range("N97").Value = -1692
$ws.Range("H102").Value = 2590.0908
$ws.Range("I102").Value = 2565.9048
$ws.Range("K102").Value = 2565.9048
$ws.Range("M102").Value = -943.9047999999998
$ws.Range("H107").Value = 351.44446
$ws.Range("J107").Value = 1049
$ws.Range("L107").Value = 1049
$ws.Range("N107").Value = -4889
$ws.Range("H122").Value = 4450.278
$ws.Range("I122").Value = 3807.7144
$ws.Range("J122").Value = 6699.25
$ws.Range("K122").Value = 11423.1432
$ws.Range("L122").Value = 20097.75
$ws.Range("M122").Value = -8973.143199999999
$ws.Range("N122").Value = -24997.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13197.2
$ws.Range("I40").Value = 12141.182
$ws.Range("K40").Value = 12141.182
$ws.Range("M40").Value = -12005.182
$ws.Range("H61").Value = 4979.8
$ws.Range("I61").Value = 4999.75
$ws.Range("J61").Value = 4900
$ws.Range("K61").Value = 4999.75
$ws.Range("L61").Value = 4900
$ws.Range("M61").Value = -4797.75
$ws.Range("N61").Value = -5304
$ws.Range("H93").Value = 2405
$ws.Range("I93").Value = 2666
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 2666
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = -1418
$ws.Range("N93").Value = -3596
$ws.Range("H113").Value = 4979.8
$ws.Range("I113").Value = 4999.75
$ws.Range("J113").Value = 4900
$ws.Range("K113").Value = 4999.75
$ws.Range("L113").Value = 4900
$ws.Range("M113").Value = -2829.75
$ws.Range("N113").Value = -9240
$ws.Range("H122").Value = 5453.75
$ws.Range("I122").Value = 5530.263
$ws.Range("K122").Value = 16590.789
$ws.Range("M122").Value = -14140.789
$ws.Range("H132").Value = 8930.362999999999
$ws.Range("I132").Value = 13666.3
$ws.Range("J132").Value = 4983.75
$ws.Range("K132").Value = 40998.89999999999
$ws.Range("L132").Value = 14951.25
$ws.Range("M132").Value = -38468.89999999999
$ws.Range("N132").Value = -20011.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1234
$ws.Range("I107").Value = 969
$ws.Range("J107").Value = 1499
$ws.Range("K107").Value = 2907
$ws.Range("L107").Value = 4497
$ws.Range("M107").Value = -987
$ws.Range("N107").Value = -8337
$ws.Range("H113").Value = 990.38464
$ws.Range("I113").Value = 1030.909
$ws.Range("J113").Value = 767.5
$ws.Range("K113").Value = 3092.727
$ws.Range("L113").Value = 2302.5
$ws.Range("M113").Value = -922.7270000000003
$ws.Range("N113").Value = -6642.5
$ws.Range("H126").Value = 1716.1786
$ws.Range("I126").Value = 1676.7778
$ws.Range("K126").Value = 5030.3334
$ws.Range("M126").Value = -2560.3334
$ws.Range("H132").Value = 6966.846
$ws.Range("I132").Value = 6395.6665
$ws.Range("K132").Value = 19186.9995
$ws.Range("M132").Value = -16656.9995

Write-Output "Applied 214 cell updates across 8 sheets"
